$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 7.5
$ws.Range("G7").Value = 1.95
$ws.Range("I7").Value = 3.6
$ws.Range("J7").Value = 2.5
$ws.Range("L7").Value = 4
$ws.Range("Q7").Value = 1.65
$ws.Range("R7").Value = 2.2
$ws.Range("Y7").Value = 8.5
$ws.Range("Z7").Value = 17
$ws.Range("AH7").Value = 21
$ws.Range("AI7").Value = 13
$ws.Range("AO7").Value = 10
$ws.Range("L8").Value = 3.25
$ws.Range("O8").Value = 1.2
$ws.Range("P8").Value = 4.33
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 2.1
$ws.Range("U8").Value = 1.57
$ws.Range("V8").Value = 2.25
$ws.Range("M13").Value = 1.13
$ws.Range("N13").Value = 6
$ws.Range("O13").Value = 1.62
$ws.Range("P13").Value = 2.2
$ws.Range("G14").Value = 1.75
$ws.Range("H14").Value = 3.75
$ws.Range("I14").Value = 4.33
$ws.Range("J14").Value = 2.3
$ws.Range("U14").Value = 1.57
$ws.Range("V14").Value = 2.25
$ws.Range("AF14").Value = 41
$ws.Range("AJ14").Value = 41
$ws.Range("AK14").Value = 29
$ws.Range("AW14").Value = 6
$ws.Range("BB14").Value = 151
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 9
$ws.Range("O17").Value = 1.4
$ws.Range("P17").Value = 2.75
$ws.Range("Q17").Value = 2.25
$ws.Range("R17").Value = 1.62
$ws.Range("O18").Value = 1.57
$ws.Range("P18").Value = 2.25
$ws.Range("U18").Value = 2.25
$ws.Range("V18").Value = 1.57
$ws.Range("AE18").Value = 21
$ws.Range("AG18").Value = 6.5
$ws.Range("BA18").Value = 126
$ws.Range("G19").Value = 1.42
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 8.5
$ws.Range("U19").Value = 2.63
$ws.Range("V19").Value = 1.44
$ws.Range("AF19").Value = 126
